$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.027.44"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  +1.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.412.93"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.19"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  +1.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.69"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +3.44%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.408.09"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("E11").Value = "  -1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.40"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.06"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +3.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  +4.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.839.12"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.893.82"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.408.29"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  +0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +3.02%  "

$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.23"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.73"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.17"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("E25").Value = "  +2.21%  "

$ws.Range("E26").Value = "  +7.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "580.56"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  +14.11%  "

$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.524.76"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  +5.37%  "

$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("E32").Value = "  +4.40%  "

$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("E34").Value = "  +2.52%  "

$ws.Range("E35").Value = "  +2.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.67"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  +5.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.76"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.89"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  +3.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.67"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  +11.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.93"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  +1.84%  "

$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0540"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  +3.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.22"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +4.55%  "

$ws.Range("E49").Value = "  +2.31%  "

$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("E51").Value = "  +1.97%  "
